# Project DesignFirst save: update cell C10 on the active sheet from 18 to 100
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 100
